# Update automàtic: dades i banners [2026-02-19 15:35]
# Refreshes the DATA_EXTRACCIO (column E) timestamps for each station row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 15:33:53"
$ws.Range("E3").Value = "2026-02-19 15:33:56"
$ws.Range("E4").Value = "2026-02-19 15:33:59"
$ws.Range("E5").Value = "2026-02-19 15:34:01"
$ws.Range("E6").Value = "2026-02-19 15:34:04"
$ws.Range("E7").Value = "2026-02-19 15:34:07"
$ws.Range("E8").Value = "2026-02-19 15:34:09"
$ws.Range("E9").Value = "2026-02-19 15:34:12"
$ws.Range("E10").Value = "2026-02-19 15:34:15"
$ws.Range("E11").Value = "2026-02-19 15:34:17"
$ws.Range("E12").Value = "2026-02-19 15:34:20"
$ws.Range("E13").Value = "2026-02-19 15:34:23"
$ws.Range("E14").Value = "2026-02-19 15:34:25"
$ws.Range("E15").Value = "2026-02-19 15:34:28"
$ws.Range("E16").Value = "2026-02-19 15:34:31"
$ws.Range("E17").Value = "2026-02-19 15:34:33"
$ws.Range("E18").Value = "2026-02-19 15:34:36"
$ws.Range("E19").Value = "2026-02-19 15:34:39"
$ws.Range("E20").Value = "2026-02-19 15:34:42"
$ws.Range("E21").Value = "2026-02-19 15:34:44"
$ws.Range("E22").Value = "2026-02-19 15:34:47"
$ws.Range("E23").Value = "2026-02-19 15:34:50"
$ws.Range("E24").Value = "2026-02-19 15:34:52"
$ws.Range("E25").Value = "2026-02-19 15:34:55"
$ws.Range("E26").Value = "2026-02-19 15:34:58"
$ws.Range("E27").Value = "2026-02-19 15:35:01"
$ws.Range("E28").Value = "2026-02-19 15:35:03"
$ws.Range("E29").Value = "2026-02-19 15:35:06"
$ws.Range("E30").Value = "2026-02-19 15:35:09"
$ws.Range("E31").Value = "2026-02-19 15:35:11"
$ws.Range("E32").Value = "2026-02-19 15:35:14"
$ws.Range("E33").Value = "2026-02-19 15:35:16"
$ws.Range("E34").Value = "2026-02-19 15:35:18"
$ws.Range("E35").Value = "2026-02-19 15:35:21"
$ws.Range("E36").Value = "2026-02-19 15:35:24"
$ws.Range("E37").Value = "2026-02-19 15:35:27"
$ws.Range("E38").Value = "2026-02-19 15:35:29"
$ws.Range("E39").Value = "2026-02-19 15:35:32"
$ws.Range("E40").Value = "2026-02-19 15:35:34"
$ws.Range("E41").Value = "2026-02-19 15:35:37"
$ws.Range("E42").Value = "2026-02-19 15:35:39"
$ws.Range("E43").Value = "2026-02-19 15:35:42"
$ws.Range("E44").Value = "2026-02-19 15:35:45"
$ws.Range("E45").Value = "2026-02-19 15:35:47"
$ws.Range("E46").Value = "2026-02-19 15:35:50"
